$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance row for 17-Aug-2023 (row 11)
$ws.Range("A11").Value = 45155
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("B11").Value = "PRESENT"
$ws.Range("C11").Value = "PRESENT"
$ws.Range("D11").Value = "PRESENT"
$ws.Range("E11").Value = "PRESENT"
$ws.Range("F11").Value = "ABSENT"
$ws.Range("G11").Value = "ABSENT"
$ws.Range("H11").Value = "ABSENT"
$ws.Range("I11").Value = "ABSENT"

# Renuka's "No Response" comments on the ABSENT cells, matching the
# other rows in the sheet.
$noResponse = "RENUKA:" + [char]10 + "No Response"
$ws.Range("F11").AddComment($noResponse)
$ws.Range("G11").AddComment($noResponse)
$ws.Range("H11").AddComment($noResponse)
$ws.Range("I11").AddComment($noResponse)

$ws.Range("I11").Select()
